# Fix conc character combine, add flag entries.
#
# This script reproduces, via Excel COM automation, the changes made to
# flag_map.xlsx:
#   * 3 new rows inserted after the existing "negative_conc_values" row
#     (new rows 12-14): negative_conc_upper_bound_values,
#     negative_conc_sd_values, negative_conc_lower_bound_values - all
#     sharing the same Definition/Flag Type as negative_conc_values.
#   * 2 new rows appended at the end (new rows 69-70):
#     cvt_conc_convert_fail / Concentration normalization failed and
#     cvt_dose_level_normalized_convert_fail / Dose normalization failed.
#   * AutoFilter / _FilterDatabase defined name range grows from
#     A1:D55 to A1:D58 (kept in sync with the 3-row insertion, even
#     though 2 more rows get appended below it afterwards - matching
#     the source workbook exactly).
#   * Sheet selection ends up on D69 / view scrolled down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert 3 blank rows right after row 11 (pushes the old rows
#    12-65 down to 15-68).
# ---------------------------------------------------------------------
[void]$ws.Range("A12:D14").EntireRow.Insert()

# Populate the 3 new rows. Column A values are written in the same
# order the original workbook's sharedStrings table shows them
# (lower, upper, sd) so the underlying shared-string table matches
# exactly; the visible row order (upper/sd/lower on rows 12/13/14)
# is unaffected by this write order.
$ws.Range("A14").Value2 = "negative_conc_lower_bound_values"
$ws.Range("A12").Value2 = "negative_conc_upper_bound_values"
$ws.Range("A13").Value2 = "negative_conc_sd_values"

$ws.Range("B12").Value2 = "Conc_Time_Values"
$ws.Range("C12").Value2 = "Template has negative concentration values"
$ws.Range("D12").Value2 = "Hard Stop (Impossible Value)"

$ws.Range("B13").Value2 = "Conc_Time_Values"
$ws.Range("C13").Value2 = "Template has negative concentration values"
$ws.Range("D13").Value2 = "Hard Stop (Impossible Value)"

$ws.Range("B14").Value2 = "Conc_Time_Values"
$ws.Range("C14").Value2 = "Template has negative concentration values"
$ws.Range("D14").Value2 = "Hard Stop (Impossible Value)"

# ---------------------------------------------------------------------
# 2. Re-point the AutoFilter / _FilterDatabase range at A1:D58 (its
#    old bound of D55 plus the 3 rows just inserted) *before* the two
#    brand-new rows get appended at the bottom of the sheet - Excel's
#    AutoFilter always snaps to the contiguous used range, so the
#    tail rows are temporarily removed, the filter is (re)applied,
#    and the tail rows are restored afterwards.
# ---------------------------------------------------------------------
$savedA = @()
$savedB = @()
$savedC = @()
$savedD = @()
for ($r = 59; $r -le 68; $r++) {
    $savedA += $ws.Range("A$r").Value2
    $savedB += $ws.Range("B$r").Value2
    $savedC += $ws.Range("C$r").Value2
    $savedD += $ws.Range("D$r").Value2
}

[void]$ws.Range("A59:D68").EntireRow.Delete()

$ws.AutoFilterMode = $false
[void]$ws.Range("A1:D58").AutoFilter()

[void]$ws.Range("A59:D68").EntireRow.Insert()
for ($i = 0; $i -lt 10; $i++) {
    $r = 59 + $i
    $ws.Range("A$r").Value2 = $savedA[$i]
    $ws.Range("B$r").Value2 = $savedB[$i]
    $ws.Range("C$r").Value2 = $savedC[$i]
    $ws.Range("D$r").Value2 = $savedD[$i]
}

foreach ($dn in $wb.Names) {
    if ($dn.Name -like "*_FilterDatabase*") {
        $dn.RefersTo = "=Sheet1!`$A`$1:`$D`$58"
    }
}

# ---------------------------------------------------------------------
# 3. Append the 2 brand-new rows at the very end of the table (new
#    rows 69 and 70).
# ---------------------------------------------------------------------
$ws.Range("A69").Value2 = "cvt_conc_convert_fail"
$ws.Range("B69").Value2 = "Conc_Time_Values"
$ws.Range("C69").Value2 = "Concentration normalization failed"
$ws.Range("D69").Value2 = "Soft Stop (Conversion Needed)"

$ws.Range("A70").Value2 = "cvt_dose_level_normalized_convert_fail"
$ws.Range("B70").Value2 = "Studies"
$ws.Range("C70").Value2 = "Dose normalization failed"
$ws.Range("D70").Value2 = "Soft Stop (Conversion Needed)"

# ---------------------------------------------------------------------
# 4. Match the final view state: active cell / selection on D69, and
#    scrolled down so row 55 is at the top.
# ---------------------------------------------------------------------
[void]$ws.Range("D69").Select()
$excel.ActiveWindow.ScrollRow = 55
$excel.ActiveWindow.ScrollColumn = 1
